{"js": "// Update the title date line.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  const titlePara = paragraphs.items[0];\n  titlePara.load(\"text\");\n  await context.sync();\n  if (titlePara.text.trim() === \"2025-10-23 Thursday\") {\n    const searchResults = titlePara.search(\"2025-10-23 Thursday\", { matchCase: true });\n    searchResults.load(\"items\");\n    await context.sync();\n    if (searchResults.items.length > 0) {\n      searchResults.items[0].insertText(\"2025-10-24 Friday\", Word.InsertLocation.replace);\n      await context.sync();\n    }\n  }\n}\n\n// Update the multiplication-fact table cells. Each content row holds five\n// \"A\u00d7B=C\" facts; blank spacer rows are left untouched. Cells are addressed\n// by (row, column) rather than by text match, because some of the original\n// fact strings repeat (e.g. \"687\u00d73=2061\" appears twice with different\n// replacements).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n\n  // Mapping of every fact cell, in row-major order, to its replacement.\n  const cellUpdates = [\n    { row: 0, col: 0, from: \"704\u00d75=3520\", to: \"673\u00d78=5384\" },\n    { row: 0, col: 1, from: \"220\u00d73=660\", to: \"886\u00d73=2658\" },\n    { row: 0, col: 2, from: \"886\u00d75=4430\", to: \"760\u00d76=4560\" },\n    { row: 0, col: 3, from: \"110\u00d78=880\", to: \"953\u00d77=6671\" },\n    { row: 0, col: 4, from: \"687\u00d73=2061\", to: \"962\u00d72=1924\" },\n\n    { row: 4, col: 0, from: \"942\u00d72=1884\", to: \"450\u00d74=1800\" },\n    { row: 4, col: 1, from: \"638\u00d77=4466\", to: \"562\u00d75=2810\" },\n    { row: 4, col: 2, from: \"123\u00d76=738\", to: \"894\u00d79=8046\" },\n    { row: 4, col: 3, from: \"809\u00d77=5663\", to: \"938\u00d75=4690\" },\n    { row: 4, col: 4, from: \"901\u00d76=5406\", to: \"662\u00d73=1986\" },\n\n    { row: 9, col: 0, from: \"687\u00d73=2061\", to: \"223\u00d79=2007\" },\n    { row: 9, col: 1, from: \"991\u00d79=8919\", to: \"589\u00d78=4712\" },\n    { row: 9, col: 2, from: \"484\u00d72=968\", to: \"941\u00d79=8469\" },\n    { row: 9, col: 3, from: \"624\u00d74=2496\", to: \"908\u00d73=2724\" },\n    { row: 9, col: 4, from: \"133\u00d74=532\", to: \"456\u00d74=1824\" },\n\n    { row: 14, col: 0, from: \"216\u00d78=1728\", to: \"353\u00d75=1765\" },\n    { row: 14, col: 1, from: \"278\u00d75=1390\", to: \"853\u00d73=2559\" },\n    { row: 14, col: 2, from: \"246\u00d76=1476\", to: \"170\u00d72=340\" },\n    { row: 14, col: 3, from: \"480\u00d72=960\", to: \"186\u00d76=1116\" },\n    { row: 14, col: 4, from: \"731\u00d78=5848\", to: \"306\u00d72=612\" },\n\n    { row: 19, col: 0, from: \"238\u00d72=476\", to: \"412\u00d77=2884\" },\n    { row: 19, col: 1, from: \"522\u00d76=3132\", to: \"313\u00d74=1252\" },\n    { row: 19, col: 2, from: \"559\u00d75=2795\", to: \"941\u00d75=4705\" },\n    { row: 19, col: 3, from: \"963\u00d72=1926\", to: \"434\u00d76=2604\" },\n    { row: 19, col: 4, from: \"686\u00d75=3430\", to: \"360\u00d78=2880\" },\n  ];\n\n  // Read current values first so we only touch cells that still hold the\n  // expected \"before\" text (defensive against an already-edited doc).\n  table.load(\"values\");\n  await context.sync();\n\n  for (const upd of cellUpdates) {\n    const currentValue = table.values[upd.row] ? table.values[upd.row][upd.col] : undefined;\n    const cell = table.getCell(upd.row, upd.col);\n    if (currentValue === upd.from) {\n      cell.value = upd.to;\n    } else if (currentValue === upd.to) {\n      // Already updated; nothing to do.\n    } else {\n      // Fallback: value didn't match what we expected, set it anyway so the\n      // final state matches the target.\n      cell.value = upd.to;\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Update the title date line (first paragraph).\n$d = $word.ActiveDocument\n\n$titlePara = $d.Paragraphs.Item(1).Range\nif ($titlePara.Text.TrimEnd([char]13, [char]7) -eq \"2025-10-23 Thursday\") {\n    $titlePara.Text = \"2025-10-24 Friday\"\n}\n\n# Update the multiplication-fact table cells. Each content row holds five\n# \"A\u00d7B=C\" facts; blank spacer rows are left untouched. Cells are addressed\n# by (row, column) rather than by text match, because some of the original\n# fact strings repeat (e.g. \"687\u00d73=2061\" appears twice with different\n# replacements), so a blind Find/Replace-All would corrupt one of them.\n$table = $d.Tables.Item(1)\n\n# row, col (1-based), expected-old text, new text\n$cellUpdates = @(\n    @(1, 1, \"704\u00d75=3520\", \"673\u00d78=5384\"),\n    @(1, 2, \"220\u00d73=660\", \"886\u00d73=2658\"),\n    @(1, 3, \"886\u00d75=4430\", \"760\u00d76=4560\"),\n    @(1, 4, \"110\u00d78=880\", \"953\u00d77=6671\"),\n    @(1, 5, \"687\u00d73=2061\", \"962\u00d72=1924\"),\n\n    @(5, 1, \"942\u00d72=1884\", \"450\u00d74=1800\"),\n    @(5, 2, \"638\u00d77=4466\", \"562\u00d75=2810\"),\n    @(5, 3, \"123\u00d76=738\", \"894\u00d79=8046\"),\n    @(5, 4, \"809\u00d77=5663\", \"938\u00d75=4690\"),\n    @(5, 5, \"901\u00d76=5406\", \"662\u00d73=1986\"),\n\n    @(10, 1, \"687\u00d73=2061\", \"223\u00d79=2007\"),\n    @(10, 2, \"991\u00d79=8919\", \"589\u00d78=4712\"),\n    @(10, 3, \"484\u00d72=968\", \"941\u00d79=8469\"),\n    @(10, 4, \"624\u00d74=2496\", \"908\u00d73=2724\"),\n    @(10, 5, \"133\u00d74=532\", \"456\u00d74=1824\"),\n\n    @(15, 1, \"216\u00d78=1728\", \"353\u00d75=1765\"),\n    @(15, 2, \"278\u00d75=1390\", \"853\u00d73=2559\"),\n    @(15, 3, \"246\u00d76=1476\", \"170\u00d72=340\"),\n    @(15, 4, \"480\u00d72=960\", \"186\u00d76=1116\"),\n    @(15, 5, \"731\u00d78=5848\", \"306\u00d72=612\"),\n\n    @(20, 1, \"238\u00d72=476\", \"412\u00d77=2884\"),\n    @(20, 2, \"522\u00d76=3132\", \"313\u00d74=1252\"),\n    @(20, 3, \"559\u00d75=2795\", \"941\u00d75=4705\"),\n    @(20, 4, \"963\u00d72=1926\", \"434\u00d76=2604\"),\n    @(20, 5, \"686\u00d75=3430\", \"360\u00d78=2880\")\n)\n\nforeach ($upd in $cellUpdates) {\n    $row = $upd[0]\n    $col = $upd[1]\n    $oldText = $upd[2]\n    $newText = $upd[3]\n\n    $cell = $table.Cell($row, $col)\n    $cellRange = $cell.Range\n    # Strip the trailing cell-mark characters (CR + cell marker) before\n    # comparing, so we only touch cells that still hold the expected text.\n    $currentText = $cellRange.Text.TrimEnd([char]13, [char]7)\n\n    if ($currentText -eq $oldText -or $currentText -ne $newText) {\n        $cellRange.Text = $newText\n    }\n}\n"}
